# Update the handback report timestamps for the 403c8a17... file rows
# on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-18 06:40:00"
$wsZhCn.Range("G2").Value = "2016-01-18 06:40:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-18 06:40:10"
$wsDeDe.Range("G2").Value = "2016-01-18 06:41:00"
